# Apply cryptos list price/volume refresh (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, preserving number-like strings
# (e.g. "1.020", "0.1530") exactly instead of letting Excel coerce them
# to numeric values, then restore the default "Normal" style so no
# stray per-cell formatting is left behind.
function Set-TextCell($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextCell 'D2' '27.242.82'
$ws.Range('E2').Value = '  +1.43%  '
Set-TextCell 'D3' '1.860.72'
$ws.Range('E3').Value = '  +1.10%  '
$ws.Range('E4').Value = '  +1.48%  '
Set-TextCell 'D5' '312.58'
$ws.Range('E5').Value = '  +1.15%  '
Set-TextCell 'D6' '1.020'
$ws.Range('E6').Value = '  +1.48%  '
Set-TextCell 'D7' '0.4792'
$ws.Range('E7').Value = '  +1.99%  '
Set-TextCell 'D8' '0.3724'
$ws.Range('E8').Value = '  +1.81%  '
Set-TextCell 'D9' '0.07328'
$ws.Range('E9').Value = '  +2.62%  '
Set-TextCell 'D10' '0.9349'
$ws.Range('E10').Value = '  +1.35%  '
Set-TextCell 'D11' '20.28'
$ws.Range('E11').Value = '  +3.80%  '
Set-TextCell 'D12' '0.07872'
$ws.Range('E12').Value = '  +2.61%  '
Set-TextCell 'D13' '1.851.51'
$ws.Range('E13').Value = '  -0.59%  '
Set-TextCell 'D14' '5.424'
$ws.Range('E14').Value = '  +2.66%  '
Set-TextCell 'D15' '6.541'
$ws.Range('E15').Value = '  +2.42%  '
Set-TextCell 'D16' '90.14'
$ws.Range('E16').Value = '  +2.24%  '
Set-TextCell 'D17' '1.023'
$ws.Range('E17').Value = '  +1.51%  '
Set-TextCell 'D18' '0.000008736'
$ws.Range('E18').Value = '  +1.27%  '
$ws.Range('E19').Value = '  +1.45%  '
Set-TextCell 'D20' '14.76'
$ws.Range('E20').Value = '  +2.22%  '
Set-TextCell 'D21' '27.277.33'
$ws.Range('E21').Value = '  +1.42%  '
Set-TextCell 'D22' '5.109'
$ws.Range('E22').Value = '  +2.13%  '
$ws.Range('E23').Value = '  +0.80%  '
Set-TextCell 'D24' '1.948'
$ws.Range('E24').Value = '  +1.35%  '
Set-TextCell 'D25' '153.82'
$ws.Range('E25').Value = '  +1.40%  '
Set-TextCell 'D26' '18.51'
$ws.Range('E26').Value = '  +1.62%  '
Set-TextCell 'D27' '2.002'
$ws.Range('E27').Value = '  -0.20%  '
Set-TextCell 'D28' '115.67'
$ws.Range('E28').Value = '  +1.43%  '
Set-TextCell 'D29' '4.986'
$ws.Range('E29').Value = '  +2.26%  '
Set-TextCell 'D30' '0.08891'
$ws.Range('E30').Value = '  +0.91%  '
Set-TextCell 'D31' '3.350'
$ws.Range('E31').Value = '  +4.39%  '
Set-TextCell 'D32' '1.188'
$ws.Range('E32').Value = '  +1.23%  '
Set-TextCell 'D33' '4.588'
$ws.Range('E33').Value = '  +2.53%  '
Set-TextCell 'D34' '0.7404'
$ws.Range('E34').Value = '  -0.61%  '
Set-TextCell 'D35' '2.678'
$ws.Range('E35').Value = '  -3.54%  '
Set-TextCell 'D36' '1.125'
$ws.Range('E36').Value = '  +3.48%  '
$ws.Range('E37').Value = '  +4.86%  '
Set-TextCell 'D38' '0.05262'
$ws.Range('E38').Value = '  +1.21%  '
Set-TextCell 'D39' '0.5345'
$ws.Range('E39').Value = '  +3.03%  '
Set-TextCell 'D40' '7.118'
$ws.Range('E40').Value = '  +2.28%  '
Set-TextCell 'D41' '0.1530'
$ws.Range('E41').Value = '  +1.51%  '
Set-TextCell 'D42' '8.330'
$ws.Range('E42').Value = '  +2.37%  '
Set-TextCell 'D43' '10.61'
$ws.Range('E43').Value = '  +1.70%  '
Set-TextCell 'D44' '0.4794'
$ws.Range('E44').Value = '  +2.35%  '
Set-TextCell 'D45' '1.021'
$ws.Range('E45').Value = '  +1.57%  '
$ws.Range('E46').Value = '  +1.41%  '
Set-TextCell 'D48' '66.42'
$ws.Range('E48').Value = '  +1.59%  '
Set-TextCell 'D49' '0.06078'
Set-TextCell 'D50' '0.9024'
$ws.Range('E50').Value = '  +1.37%  '
Set-TextCell 'D51' '36.70'
$ws.Range('E51').Value = '  +1.58%  '
